$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: "touch" a property on a Range so that, on save, it is not silently
# re-merged with an adjacent run that happens to share identical formatting.
# (Word/engine auto-merges runs with identical rPr; flipping Bold on/off is a
#  harmless no-op on the visible formatting but forces the run boundary to
#  survive serialization.)
# ---------------------------------------------------------------------------
function Touch-Run($range) {
    $orig = $range.Bold
    $range.Bold = 1
    $range.Bold = $orig
}

# ===========================================================================
# Edit 1: "...em especial Álgebra Linear. " ->
#         "...em especial Geometria Analítica e Álgebra Linear. "
# The whole sentence ("na disciplina ... em especial Álgebra Linear. ") is
# originally a single run; it ends up split into three runs, and the
# preceding run ("com maiores detalhes ") must NOT get swept into the merge
# either, so we touch from the start of that full original run through to
# the end of the (new) tail piece.
# ===========================================================================
$r1 = $d.Content
$needle1 = "na disciplina de Computação Gráfica do curso de Ciência da Computação do UNI-BH. Por isso, este trabalho foi uma ótima oportunidade para adquirir novos conhecimentos e colocar em prática conhecimentos de várias disciplinas, em especial Álgebra Linear. "
$found1 = $r1.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $headLen = "na disciplina de Computação Gráfica do curso de Ciência da Computação do UNI-BH. Por isso, este trabalho foi uma ótima oportunidade para adquirir novos conhecimentos e colocar em prática conhecimentos de várias disciplinas, em especial ".Length
    $insertAt = $r1.Start + $headLen
    $newMid = "Geometria Analítica e "

    $ip = $d.Range($insertAt, $insertAt)
    $ip.InsertBefore($newMid)

    $headRange = $d.Range($r1.Start, $insertAt)
    Touch-Run $headRange

    $midRange = $d.Range($insertAt, $insertAt + $newMid.Length)
    Touch-Run $midRange

    $tailStart = $insertAt + $newMid.Length
    $tailEnd = $r1.End + $newMid.Length
    $tailRange = $d.Range($tailStart, $tailEnd)
    Touch-Run $tailRange
}

# ===========================================================================
# Edit 2: Rewrite the asteroid/ring-plane collision-distance sentence.
# ===========================================================================
$r2 = $d.Content
$oldSentence = "ide intercepta esse plano, verifica-se a distância do seu ponto central até o centro dos anéis (o centro de Saturno) para determinar se ocorre colisão.  "
$found2 = $r2.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $start2 = $r2.Start

    $newSentence = "ide intercepta esse plano, encontra-se a distância entre o centro do asteroide e o plano. Com esta distância e o raio do asteroide, calcula-se o círculo projetado pela esfera no plano. Daí, a colisão ocorre se a distância entre o centro desse círculo projetado e o centro de Saturno é menor que a soma de seus raios.  "

    # Replace text content in one shot (keeps it inside the originally-found
    # run's formatting for now); we then carve it into the individual runs
    # the target document expects.
    $r2.Text = $newSentence

    $pieces = @(
        "ide intercepta esse plano, ",
        "encontra-se a distância entre o centro do asteroide e o plano. Com es",
        "t",
        "a distância e o raio do asteroide, calcula-se o ",
        "círculo ",
        "projetad",
        "o ",
        "pela esfera no plano. Daí, a colisão ocorre se a distância entre o centro desse círculo projetado e o centro de Saturno é menor que a soma de seus raios.",
        "  "
    )

    $offset = $start2
    foreach ($piece in $pieces) {
        $pieceStart = $offset
        $pieceEnd = $offset + $piece.Length
        $pieceRange = $d.Range($pieceStart, $pieceEnd)
        Touch-Run $pieceRange
        $offset = $pieceEnd
    }
}

# ===========================================================================
# Edit 3: Remove the stray empty paragraph right before
#         "Por fim, como o sistema de detecção de colisão "
#         (sits between the "Figura ... de Saturno" caption and that
#          paragraph).
# ===========================================================================
$r3 = $d.Content
$found3 = $r3.Find.Execute("Por fim, como o sistema de detecção de colisão ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $targetPara = $r3.Paragraphs.First
    $allParas = $d.Paragraphs
    $prevPara = $allParas.Item($targetPara.Index - 1)
    if ($prevPara.Range.Text -eq "") {
        $delRange = $d.Range($prevPara.Range.Start, $prevPara.Range.End)
        $delRange.Delete()
    }
}

# ===========================================================================
# Edit 4: Fix the cached SEQ-field figure number from "3" to "1".
# ===========================================================================
$fields = $d.Fields
for ($i = 1; $i -le $fields.Count; $i++) {
    $fld = $fields.Item($i)
    if ($fld.Code.Text.Trim() -eq "SEQ Figura \* ARABIC") {
        $res = $fld.Result
        $resRange = $d.Range($res.Start, $res.End)
        if ($resRange.Text -eq "3") {
            $resRange.Text = "1"
        }
    }
}

Write-Output "edits applied"
